# Updates cryptocurrency Price (D) and Volume(1h) (E) columns on Sheet1
# to match the latest scraped snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.003.09"
$ws.Range("E2").Value = "  -2.78%  "
$ws.Range("D3").Value = "3.028.94"
$ws.Range("E3").Value = "  -2.12%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.21"
$ws.Range("E5").Value = "  +2.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.47"
$ws.Range("E6").Value = "  -0.81%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.025.21"
$ws.Range("E8").Value = "  -2.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.498"
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.150"
$ws.Range("E10").Value = "  -4.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.09"
$ws.Range("E11").Value = "  -5.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.452"
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000223"
$ws.Range("E13").Value = "  -1.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.55"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "3.514.41"
$ws.Range("E15").Value = "  -2.08%  "
$ws.Range("D16").Value = "62.081.94"
$ws.Range("E16").Value = "  -2.67%  "
$ws.Range("E17").Value = "  -2.45%  "
$ws.Range("D18").Value = "3.029.73"
$ws.Range("E18").Value = "  -1.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.69"
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "475.45"
$ws.Range("E20").Value = "  -1.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.33"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.678"
$ws.Range("E22").Value = "  -3.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.10"
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.52"
$ws.Range("E24").Value = "  +1.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.18"
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.75"
$ws.Range("E27").Value = "  +1.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.84"
$ws.Range("E28").Value = "  -3.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.92"
$ws.Range("E30").Value = "  +1.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.87"
$ws.Range("E31").Value = "  -1.55%  "
$ws.Range("E32").Value = "  +1.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.33"
$ws.Range("E33").Value = "  -0.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "55.90"
$ws.Range("E34").Value = "  -2.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.48"
$ws.Range("E35").Value = "  +2.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.94"
$ws.Range("E36").Value = "  -1.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "460.54"
$ws.Range("E37").Value = "  -7.54%  "
$ws.Range("D38").Value = "3.215.15"
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0801"
$ws.Range("E39").Value = "  +0.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0387"
$ws.Range("E40").Value = "  -4.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.119"
$ws.Range("E41").Value = "  +0.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.17"
$ws.Range("E42").Value = "  +0.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.48"
$ws.Range("E43").Value = "  -6.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "26.00"
$ws.Range("E45").Value = "  +5.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.246"
$ws.Range("E46").Value = "  -3.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.01"
$ws.Range("E47").Value = "  -1.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.109"
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "118.69"
$ws.Range("E49").Value = "  -3.72%  "
$ws.Range("D50").Value = "0.0₃0499"
$ws.Range("E50").Value = "  -5.70%  "
$ws.Range("E51").Value = "  +7.18%  "
